$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row for the renamed taxon "Eggerthella sp. type 1" and
# populate only the current_classification / new_classification columns
# (no tax_id / rank values are known for this entry).
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "Eggerthella sp. type 1 (species)"
$ws.Range("D6").Value = "Eggerthella sp. type 1"

# Re-sort the data rows (A2:D14) alphabetically by current_classification,
# which is where the new row belongs anyway and is what leaves the
# worksheet's sortState behind.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A14")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:D14"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Formatting got normalized back to the default style during the edit.
$ws.Range("A1:D14").ClearFormats()

# Leave the selection where the user's cursor ended up.
$ws.Range("C18").Select() | Out-Null
